$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = 1.91
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 4.5
$ws.Range("J6").Value = 2.75
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5
$ws.Range("X6").Value = 8
$ws.Range("AN6").Value = 3.75
$ws.Range("AX6").Value = 26
